$d = $word.ActiveDocument

# Locate the "Requisitos" heading paragraph, then operate on the bullet-list
# paragraph right after it (the one containing all the "XXXXNNNN - ... (Requisito)"
# lines). Finding it by heading text is more robust than a hard-coded paragraph index.
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd([char]13) -eq "Requisitos") {
        $target = $p.Next(1)
        break
    }
}

if ($target -eq $null) {
    throw "Could not locate the Requisitos list paragraph"
}

$full = $target.Range
$start = $full.Start
$end = $full.End

# Range covering the paragraph's content but NOT its trailing paragraph mark,
# so InsertXML only replaces the bullet text/runs and leaves the paragraph
# (and its ListBullet style) intact.
$contentRange = $d.Range($start, $end - 1)

$xml = @"
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p><w:r><w:t>LOQ4100 -  Fundamentos de Química para Engenharia I (Requisito)</w:t><w:br/></w:r><w:r><w:t>LOQ4098 -  Fundamentos de Química para Engenharia II (Requisito)</w:t><w:br/></w:r><w:r><w:t>LOB1053 -  Física III  (Requisito)</w:t><w:br/></w:r><w:r><w:t>LOB1038 -  Física Experimental I  (Requisito)</w:t><w:br/></w:r><w:r><w:t>LOM3013 -  Ciência dos Materiais  (Requisito)</w:t><w:br/></w:r><w:r><w:t>LOB1018 -  Física I  (Requisito)</w:t><w:br/></w:r><w:r><w:t>LOB1041 -  Física Experimental II  (Requisito)</w:t><w:br/></w:r><w:r><w:t>LOQ4095 -  Química Geral Experimental  (Requisito)</w:t><w:br/></w:r><w:r><w:t>LOB1036 -  Geometria Analítica  (Requisito)</w:t><w:br/></w:r><w:r><w:t>LOM3037 -  Química Inorgânica  (Requisito)</w:t><w:br/></w:r><w:r><w:t>LOM3099 -  Estática  (Requisito)</w:t><w:br/></w:r><w:r><w:t>LOM3204 -  Desenho Técnico e Projeto Assistido por Computador  (Requisito)</w:t><w:br/></w:r><w:r><w:t>LOM3056 -  Fundamentos de Química Orgânica  (Requisito)</w:t><w:br/></w:r><w:r><w:t>LOQ4246 -  Engenharia da Sustentabilidade  (Requisito)</w:t><w:br/></w:r><w:r><w:t>LOM3018 -  Introdução à Engenharia de Materiais  (Requisito)</w:t><w:br/></w:r><w:r><w:t>LOB1012 -  Estatística  (Requisito)</w:t><w:br/></w:r><w:r><w:t>LOB1004 -  Cálculo II  (Requisito)</w:t><w:br/></w:r><w:r><w:t>LOB1046 -  Engenharia do Meio Ambiente  (Requisito)</w:t><w:br/></w:r><w:r><w:t>LOB1003 -  Cálculo I  (Requisito)</w:t><w:br/></w:r><w:r><w:t>LOB1045 -  Leitura e Produção de Textos Acadêmicos  (Requisito)</w:t><w:br/></w:r><w:r><w:t>LOB1052 -  Cálculo III  (Requisito)</w:t><w:br/></w:r><w:r><w:t>LOM3104 -  Projeto Integrado em Engenharia de Materiais I  (Requisito)</w:t><w:br/></w:r><w:r><w:t>LOB1008 -  Ciência, Tecnologia e Sociedade  (Requisito)</w:t><w:br/></w:r><w:r><w:t>LOB1037 -  Àlgebra Linear  (Requisito)</w:t><w:br/></w:r><w:r><w:t>LOB1019 -  Física II  (Requisito)</w:t><w:br/></w:r><w:r><w:t>LOB1039 -  Física Experimental III  (Requisito)</w:t><w:br/></w:r><w:r><w:t>LOM3105 -  Computação e análise de dados em Engenharia  (Requisito)</w:t><w:br/></w:r></w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
"@

$contentRange.InsertXML($xml)
